# Apply the BOM refresh: remove the obsolete "7513D2-L" (LightPipe) line item,
# which shifts every subsequent row up by one, and refresh the "Quantity
# Available" (column K) stock figures across the sheet. Also narrow the
# Manufacturer column (B) to fit the new content.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Narrow the "Manufacturer" column. Excel's ColumnWidth (character units)
# is offset from the OOXML <col width> by ~5/6 of a character, so subtract
# that to land on an exact width of 27 in the saved file.
$ws.Columns.Item(2).ColumnWidth = 26.16666666666667

# The LightPipe (7513D2-L) component was dropped from the design; deleting its
# row shifts rows 14:25 up into 13:24 automatically.
$ws.Rows.Item(13).Delete()

# Refresh the "Quantity Available" column (K) with the latest stock figures.
$qtyAvailable = @{
    2  = 1362
    3  = 1013
    4  = 1772
    5  = 5706
    6  = 115422
    7  = 181978
    8  = 16344
    9  = 106442
    10 = 3139
    11 = 25643
    12 = 17756
    13 = 15875
    14 = 337750
    15 = 5562
    16 = 14785
    17 = 4495
    18 = 3663
    19 = 791941
    20 = 1333942
    21 = 1139
    22 = 199128
    23 = 3384
    24 = 676
}

foreach ($row in $qtyAvailable.Keys) {
    $ws.Cells.Item($row, 11).Value = $qtyAvailable[$row]
}
